$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (values like "1.001", "311.01", etc.)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update coin price (column D) and volume/change (column E) values
$ws.Range("D2").Value = '23.868.52'
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").Value = '1.652.92'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '311.01'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.3817'
$ws.Range("E8").Value = '  -2.86%  '
$ws.Range("D9").Value = '51.59'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("D10").Value = '1.347'
$ws.Range("E10").Value = '  -3.53%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '0.08477'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '23.99'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").Value = '7.055'
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").Value = '8.065'
$ws.Range("E15").Value = '  +1.56%  '
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '1.664.60'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '94.18'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = '0.06990'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '19.62'
$ws.Range("E20").Value = '  -4.93%  '
$ws.Range("D21").Value = '6.995'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '13.70'
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").Value = '23.880.35'
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("D25").Value = '2.432'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").Value = '2.969'
$ws.Range("E26").Value = '  -3.25%  '
$ws.Range("D27").Value = '22.06'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '153.99'
$ws.Range("E28").Value = '  -2.21%  '
$ws.Range("D29").Value = '5.431'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").Value = '137.91'
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("D31").Value = '7.806'
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").Value = '2.485'
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("D33").Value = '1.837.77'
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = '0.08174'
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").Value = '1.010'
$ws.Range("E35").Value = '  -5.49%  '
$ws.Range("D36").Value = '0.02914'
$ws.Range("E36").Value = '  -5.94%  '
$ws.Range("D37").Value = '6.665'
$ws.Range("E37").Value = '  -3.74%  '
$ws.Range("D38").Value = '10.78'
$ws.Range("E38").Value = '  -3.45%  '
$ws.Range("D39").Value = '0.2676'
$ws.Range("E39").Value = '  -3.35%  '
$ws.Range("D40").Value = '0.09139'
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '13.56'
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.7565'
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").Value = '16.65'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '0.6939'
$ws.Range("E45").Value = '  -2.69%  '
$ws.Range("D46").Value = '2.452'
$ws.Range("E46").Value = '  -4.03%  '
$ws.Range("D47").Value = '4.105'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").Value = '0.9987'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").Value = '0.08298'
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").Value = '133.30'
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("D51").Value = '1.229'
$ws.Range("E51").Value = '  -3.42%  '
